# Gallery_SlotCards_ProPanels.xlsx
# "Added Pro lite test data sheet and updated other test data related to pro lite panels"
#
# 1. Remove the stray "Sheet2" tab (duplicate of the Germany slot-cards sheet).
# 2. Rename "Slot Cards 215 Panel" -> "Germany_SlotCards215Panel" and update its
#    NGC story text to note it also applies to the ProLite panel.
# 3. Duplicate that sheet into a new "Czech_SlotCards215Panel" tab for the Czech
#    market, with its own market label + NGC story text.
# 4. Tidy up the selections / active tab so the new sheet ends up active, matching
#    how Excel leaves things after these interactive edits.

$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

# --- 1. drop the redundant "Sheet2" tab -------------------------------------
$wb.Worksheets.Item("Sheet2").Delete() | Out-Null

# --- 2. rename + update the Germany slot-cards panel sheet ------------------
$germanyPanel = $wb.Worksheets.Item("Slot Cards 215 Panel")
$germanyPanel.Name = "Germany_SlotCards215Panel"
$germanyPanel.Range("B4").Value = "NGC-3475/T1730 and applicable to ProLite Panel as well"

# --- 3. duplicate it to make the Czech slot-cards panel sheet ---------------
$germanyPanel.Copy($null, $germanyPanel) | Out-Null
$czechPanel = $wb.Worksheets.Item($germanyPanel.Index + 1)
$czechPanel.Name = "Czech_SlotCards215Panel"
$czechPanel.Range("B2").Value = "Czech Market"
$czechPanel.Range("B4").Value = "NGC-3477/T1734 and applicable to ProLite Panel as well"

# --- 4. refresh selections so the UI state matches the post-edit workbook ---
$czechSheet = $wb.Worksheets.Item("Czech")
$czechSheet.Activate()
$czechSheet.Range("B2").Select() | Out-Null

$czechPanel.Activate()
$czechPanel.Range("A9").Select() | Out-Null
